$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 46, shifting existing rows 46-110 down to 47-111
$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 44540
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100103
$ws.Cells.Item(46, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(46, 9).Value = 100103004
$ws.Cells.Item(46, 10).Value = "Durazno"
$ws.Cells.Item(46, 11).Value = "Royal Glory"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 220
$ws.Cells.Item(46, 14).Value = 12000
$ws.Cells.Item(46, 15).Value = 13000
$ws.Cells.Item(46, 16).Value = 12455
$ws.Cells.Item(46, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 830
$ws.Cells.Item(46, 20).Value = 15
